$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute('Umbizo limesahihishwa sio wakati', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'Format has been corrected not the timing'
} else {
    Write-Host 'NOT FOUND (#1): Umbizo limesahihishwa sio wakati'
}

$rng = $d.Content
$found = $rng.Find.Execute('Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'I added 25 seconds to each timing to correct for the intro song -john argentino'
} else {
    Write-Host 'NOT FOUND (#2): Niliongeza sekunde 25 kwa kila muda ili kusahihisha wimbo wa utangulizi -john argentino'
}

$rng = $d.Content
$found = $rng.Find.Execute('Tatizo la uwanja wa ndege - manukuu:', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'The airport problem - subtitles:'
} else {
    Write-Host 'NOT FOUND (#3): Tatizo la uwanja wa ndege - manukuu:'
}

$rng = $d.Content
$found = $rng.Find.Execute('Utawala wa tatu', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'The administrations of three'
} else {
    Write-Host 'NOT FOUND (#4): Utawala wa tatu'
}

$rng = $d.Content
$found = $rng.Find.Execute('miji jirani: A, B na C waliamua', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'neighboring cities: A, B and C decided'
} else {
    Write-Host 'NOT FOUND (#5): miji jirani: A, B na C waliamua'
}

$rng = $d.Content
$found = $rng.Find.Execute('kujenga uwanja wa ndege unaogawanya gharama za', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'to build an airport dividing the costs of'
} else {
    Write-Host 'NOT FOUND (#6): kujenga uwanja wa ndege unaogawanya gharama za'
}

$rng = $d.Content
$found = $rng.Find.Execute('utekelezaji. Hali juu ya', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'implementation. The condition on the'
} else {
    Write-Host 'NOT FOUND (#7): utekelezaji. Hali juu ya'
}

$rng = $d.Content
$found = $rng.Find.Execute('uchaguzi wa mahali pa kufaa zaidi ni', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'choice of the most suitable place is'
} else {
    Write-Host 'NOT FOUND (#8): uchaguzi wa mahali pa kufaa zaidi ni'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwamba jumla ya umbali kutoka kwa kila mmoja', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'that the sum of the distances from each'
} else {
    Write-Host 'NOT FOUND (#9): kwamba jumla ya umbali kutoka kwa kila mmoja'
}

$rng = $d.Content
$found = $rng.Find.Execute('mji kwa uwanja wa ndege ni ndogo kama', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'city to the airport is as small as'
} else {
    Write-Host 'NOT FOUND (#10): mji kwa uwanja wa ndege ni ndogo kama'
}

$rng = $d.Content
$found = $rng.Find.Execute('inawezekana. Timu ya wataalam wanaohusika', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'possible. The team of experts in charge'
} else {
    Write-Host 'NOT FOUND (#11): inawezekana. Timu ya wataalam wanaohusika'
}

$rng = $d.Content
$found = $rng.Find.Execute('ya kazi imeunda mfano wa kupata', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'of the work has created a model to get'
} else {
    Write-Host 'NOT FOUND (#12): ya kazi imeunda mfano wa kupata'
}

$rng = $d.Content
$found = $rng.Find.Execute('wazo la awali la mahali pa kuweka', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'a preliminary idea of where to place the'
} else {
    Write-Host 'NOT FOUND (#13): wazo la awali la mahali pa kuweka'
}

$rng = $d.Content
$found = $rng.Find.Execute('muundo. Ovyo wao wapo', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'structure. At their disposal there are'
} else {
    Write-Host 'NOT FOUND (#14): muundo. Ovyo wao wapo'
}

$rng = $d.Content
$found = $rng.Find.Execute('konokono wengine pete kubwa ya chuma na ndefu', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'some snails a big metal ring and a long'
} else {
    Write-Host 'NOT FOUND (#15): konokono wengine pete kubwa ya chuma na ndefu'
}

$rng = $d.Content
$found = $rng.Find.Execute('kamba.', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'string.'
} else {
    Write-Host 'NOT FOUND (#16): kamba.'
}

$rng = $d.Content
$found = $rng.Find.Execute('Eleza jinsi timu inaweza kusimamia matumizi', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'Explain how the team can manage to use'
} else {
    Write-Host 'NOT FOUND (#17): Eleza jinsi timu inaweza kusimamia matumizi'
}

$rng = $d.Content
$found = $rng.Find.Execute('nyenzo za kusema takriban', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'the materials to tell approximately the'
} else {
    Write-Host 'NOT FOUND (#18): nyenzo za kusema takriban'
}

$rng = $d.Content
$found = $rng.Find.Execute('eneo bora la uwanja wa ndege. Fikiria', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'ideal location of the airport. Imagine'
} else {
    Write-Host 'NOT FOUND (#19): eneo bora la uwanja wa ndege. Fikiria'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwamba miji imewekwa kwenye', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'that the cities are placed at the'
} else {
    Write-Host 'NOT FOUND (#20): kwamba miji imewekwa kwenye'
}

$rng = $d.Content
$found = $rng.Find.Execute('vipeo vya pembetatu ambayo ni', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'vertices of a triangle which is'
} else {
    Write-Host 'NOT FOUND (#21): vipeo vya pembetatu ambayo ni'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwa hakika imetolewa tena kwa kiwango kama', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'obviously reproduced in scale as'
} else {
    Write-Host 'NOT FOUND (#22): kwa hakika imetolewa tena kwa kiwango kama'
}

$rng = $d.Content
$found = $rng.Find.Execute('inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'shown in figure. This is one possible'
} else {
    Write-Host 'NOT FOUND (#23): inavyoonyeshwa kwenye takwimu. Hili ni moja linalowezekana'
}

$rng = $d.Content
$found = $rng.Find.Execute('kuweka kamba huanza kutoka msumari mmoja,', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'setting the rope starts from one nail,'
} else {
    Write-Host 'NOT FOUND (#24): kuweka kamba huanza kutoka msumari mmoja,'
}

$rng = $d.Content
$found = $rng.Find.Execute('huenda ndani ya pete, huzunguka', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'goes inside the ring, goes around the'
} else {
    Write-Host 'NOT FOUND (#25): huenda ndani ya pete, huzunguka'
}

$rng = $d.Content
$found = $rng.Find.Execute('msumari mwingine, msumari wa tatu, ndani ya', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'other nail, the third nail, inside the'
} else {
    Write-Host 'NOT FOUND (#26): msumari mwingine, msumari wa tatu, ndani ya'
}

$rng = $d.Content
$found = $rng.Find.Execute('pete tena na sasa unaweza kuvuta tu', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'ring again and now you can just pull the'
} else {
    Write-Host 'NOT FOUND (#27): pete tena na sasa unaweza kuvuta tu'
}

$rng = $d.Content
$found = $rng.Find.Execute('kamba ili kupata uhakika huo', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'rope in order to find the point that'
} else {
    Write-Host 'NOT FOUND (#28): kamba ili kupata uhakika huo'
}

$rng = $d.Content
$found = $rng.Find.Execute('unatafuta. Ili kufikia', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'you''re looking for. In order to reach the'
} else {
    Write-Host 'NOT FOUND (#29): unatafuta. Ili kufikia'
}

$rng = $d.Content
$found = $rng.Find.Execute('uhakika, tunapaswa kusonga kamba kidogo', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'point, we have to move the rope a bit'
} else {
    Write-Host 'NOT FOUND (#30): uhakika, tunapaswa kusonga kamba kidogo'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwa nyenzo ambazo tunatumia lakini', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'by the materials that we are using but'
} else {
    Write-Host 'NOT FOUND (#31): kwa nyenzo ambazo tunatumia lakini'
}

$rng = $d.Content
$found = $rng.Find.Execute('baada ya muda utafikia nafasi kutoka', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'after a while you''ll reach a position from'
} else {
    Write-Host 'NOT FOUND (#32): baada ya muda utafikia nafasi kutoka'
}

$rng = $d.Content
$found = $rng.Find.Execute('ambayo pete haisogei tena,', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'which the ring doesn''t move anymore,'
} else {
    Write-Host 'NOT FOUND (#33): ambayo pete haisogei tena,'
}

$rng = $d.Content
$found = $rng.Find.Execute('ambayo ni zaidi au chini ya hii. Na kama', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'which is more or less this one. And as'
} else {
    Write-Host 'NOT FOUND (#34): ambayo ni zaidi au chini ya hii. Na kama'
}

$rng = $d.Content
$found = $rng.Find.Execute('kati ya pete na misumari ni', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'between the ring and the nails are'
} else {
    Write-Host 'NOT FOUND (#35): kati ya pete na misumari ni'
}

$rng = $d.Content
$found = $rng.Find.Execute('kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'placed more or less 120 degrees from one'
} else {
    Write-Host 'NOT FOUND (#36): kuwekwa zaidi au chini ya digrii 120 kutoka kwa moja'
}

$rng = $d.Content
$found = $rng.Find.Execute('nyingine ambayo ni 1/3 ya mduara,', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'another which is 1/3 of a circumference,'
} else {
    Write-Host 'NOT FOUND (#37): nyingine ambayo ni 1/3 ya mduara,'
}

$rng = $d.Content
$found = $rng.Find.Execute('na hiyo ndiyo hatua tunayoiangalia', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'and that''s the point that we''re looking'
} else {
    Write-Host 'NOT FOUND (#38): na hiyo ndiyo hatua tunayoiangalia'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwa: umbali wa chini kati ya', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'for: the minimum distance between the'
} else {
    Write-Host 'NOT FOUND (#39): kwa: umbali wa chini kati ya'
}

$rng = $d.Content
$found = $rng.Find.Execute('misumari na uwanja wa ndege unapojumlisha', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'nails and the airport when you sum it'
} else {
    Write-Host 'NOT FOUND (#40): misumari na uwanja wa ndege unapojumlisha'
}

$rng = $d.Content
$found = $rng.Find.Execute('kwa sababu kuna ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'because there is some '
} else {
    Write-Host 'NOT FOUND (#41): kwa sababu kuna '
}

$rng = $d.Content
$found = $rng.Find.Execute('upinzani', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'resistance'
} else {
    Write-Host 'NOT FOUND (#42): upinzani'
}

$rng = $d.Content
$found = $rng.Find.Execute(' uliosababishwa', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = ' caused'
} else {
    Write-Host 'NOT FOUND (#43):  uliosababishwa'
}

$rng = $d.Content
$found = $rng.Find.Execute('pamoja', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = 'ogether'
} else {
    Write-Host 'NOT FOUND (#44): pamoja'
}

$rng = $d.Content
$found = $rng.Find.Execute('[Muziki]', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $rng.Text = '[Music]'
} else {
    Write-Host 'NOT FOUND (#45): [Muziki]'
}

Write-Host "All replacements applied."
